$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.793.59"
$ws.Range("E2").Value = "  +1.72%  "

$ws.Range("D3").Value = "1.880.20"
$ws.Range("E3").Value = "  +1.27%  "

$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "332.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.53%  "

$ws.Range("E6").Value = "  +0.08%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4727"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.46%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3959"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.34%  "

$ws.Range("E9").Value = "  -1.06%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08068"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.66%  "

$ws.Range("E11").Value = "  +1.78%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.25"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.01%  "

$ws.Range("D13").Value = "1.873.97"
$ws.Range("E13").Value = "  +1.50%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.983"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.29%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.144"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.27%  "

$ws.Range("E16").Value = "  +0.42%  "

$ws.Range("E17").Value = "  +2.09%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "87.31"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.69%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06674"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.78%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.33"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.19%  "

$ws.Range("E21").Value = "  +0.12%  "

$ws.Range("D22").Value = "27.802.42"
$ws.Range("E22").Value = "  +1.79%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.545"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.13%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.01"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.62%  "

$ws.Range("E25").Value = "  +0.77%  "

$ws.Range("D26").Value = "2.101.01"
$ws.Range("E26").Value = "  +1.45%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "159.54"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.86%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.28"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.86%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.111"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.06%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.619"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.94%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "122.33"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.05%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9863"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.27%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09552"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.69%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.453"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.60%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.598"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.02%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.383"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.09%  "

$ws.Range("E37").Value = "  +2.37%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02261"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.71%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.239"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.95%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.163"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.62%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6043"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.23%  "

$ws.Range("E42").Value = "  +0.09%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1907"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.88%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.33"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.16%  "

$ws.Range("B45").Value = "WEMIXTOKEN"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.261"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.84%  "

$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5749"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.22%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.30"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.36%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.953"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.82%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.379"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.19%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06914"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.35%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "114.24"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.96%  "
